# Weekly price update: a new observation is inserted as the most recent
# row for this market/category subset. All existing rows shift down by
# one, so insert a fresh row at 167 and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("167:167").Insert()

$ws.Range("A167").Value = 8
$ws.Range("B167").Value = "Terminal La Palmera de La Serena"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 44795
$ws.Range("E167").Value = 4
$ws.Range("F167").Value = 100112037
$ws.Range("G167").Value = "Cebollín"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 2000
$ws.Range("K167").Value = 1400
$ws.Range("L167").Value = 1600
$ws.Range("M167").Value = 1500
$ws.Range("N167").Value = "`$/paquete 6 unidades"
$ws.Range("O167").Value = "Provincia del Elquí"
$ws.Range("P167").Value = 250
$ws.Range("Q167").Value = 6
$ws.Range("R167").Value = "Hortaliza"
